# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table switches from the custom "Table_0" style to the
#    built-in table style {F5AE3A67-A773-44C5-B457-EBF54144C05D}.
# 2) The presentation's theme palette is swapped from the "Integral /
#    Red Violet" palette to the default "Office Theme / Office" palette
#    (the two a:clrScheme color values that used to live in theme1.xml /
#    theme2.xml trade places).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F5AE3A67-A773-44C5-B457-EBF54144C05D}")
    }
}

# --- 2. Swap the theme color scheme ----------------------------------------
# Before: dk1 000000 lt1 FFFFFF dk2 454551 lt2 D8D9DC accent1 E32D91
#         accent2 C830CC accent3 4EA6DC accent4 4775E7 accent5 8971E1
#         accent6 D54773 hlink 6B9F25 folHlink 8C8C8C   (Integral / Red Violet)
# After:  dk1 000000 lt1 FFFFFF dk2 44546A lt2 E7E6E6 accent1 5B9BD5
#         accent2 ED7D31 accent3 A5A5A5 accent4 FFC000 accent5 4472C4
#         accent6 70AD47 hlink 0563C1 folHlink 954F72  (Office Theme / Office)
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72

# Try to rename the theme / color scheme to match (best-effort; some hosts
# keep these read-only, but attempt it so the names line up if supported).
$theme.Name = "Office Theme"
$colors.Name = "Office"
$p.Designs.Item(1).Name = "Office Theme"
